$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2;  B = -0.04915694638892727; C = 0.6586326132859144; D = 1.297461883598808;  E = 1.13906184362343;   F = 1.169184523175352;  G = 19 },
    @{ Row = 3;  B = -0.07072591166527165; C = 0.5579505877006946; D = 0.6529064599314954;  E = 0.8080262742828945; F = 0.8282610236588537; G = 18 },
    @{ Row = 4;  B = 0.03916263344646955;  C = 0.552040872335158;  D = 0.6357226611696118;  E = 0.7973221815361792; F = 0.8208689046632819; G = 17 },
    @{ Row = 5;  B = 0.05955592157360488;  C = 0.532271854806013;  D = 0.7198117377579547;  E = 0.8484171955812509; F = 0.8740799841134135; G = 16 },
    @{ Row = 6;  B = 0.09235442246305205;  C = 0.5824508295849771; D = 0.7506227629998017;  E = 0.8663848815623468; F = 0.8916838547205332; G = 15 },
    @{ Row = 7;  B = 0.1521483464337763;   C = 0.6635260061488867; D = 1.004366658952618;   E = 1.002180951202236;  F = 1.0279571060649;    G = 14 },
    @{ Row = 8;  B = 0.159095420907867;    C = 0.6746913222580827; D = 0.8632158529041898;  E = 0.9290941033631576; F = 0.9527485895797485; G = 13 },
    @{ Row = 9;  B = 0.2207410413997515;   C = 0.686169913841685;  D = 0.7776458313225046;  E = 0.8818422939066285; F = 0.8917312420012032; G = 12 },
    @{ Row = 10; B = 0.220829446489595;    C = 0.7371449325608843; D = 1.012670326230143;   E = 1.006315222099985;  F = 1.029706341620059;  G = 11 },
    @{ Row = 11; B = 0.1936765229212045;   C = 0.6934653330366425; D = 1.055080256735366;   E = 1.027170996833227;  F = 1.063312257473744;  G = 10 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
}
